$wb = $excel.ActiveWorkbook

# Sheet ALC, row 9 (G9=5487)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 93.333336
$ws.Range("I9").Value = 99.09090999999999
$ws.Range("K9").Value = 99.09090999999999
$ws.Range("M9").Value = 69.90909000000001

# Sheet ALC, row 33 (G33=5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1697.5714
$ws.Range("I33").Value = 672.4706
$ws.Range("K33").Value = 672.4706
$ws.Range("M33").Value = -443.4706

# Sheet ALC, row 113 (G113=27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7230.3076
$ws.Range("I113").Value = 5132.3335
$ws.Range("K113").Value = 5132.3335
$ws.Range("M113").Value = -1878.3335

# Sheet ARM, row 2 (G2=27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6864.9
$ws.Range("I2").Value = 5223.5884
$ws.Range("K2").Value = 5223.5884
$ws.Range("M2").Value = -5110.5884

# Sheet ARM, row 32 (G32=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1966.7273
$ws.Range("I32").Value = 1529.129
$ws.Range("J32").Value = 8749.5
$ws.Range("K32").Value = 1529.129
$ws.Range("L32").Value = 8749.5
$ws.Range("M32").Value = -1242.129
$ws.Range("N32").Value = -9323.5

# Sheet ARM, row 38 (G38=2260)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 4650
$ws.Range("I38").Value = 4650
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 4650
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -4183
$ws.Range("N38").ClearContents()

# Sheet ARM, row 74 (G74=44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4160.927
$ws.Range("I74").Value = 4017.6365
$ws.Range("J74").Value = 4326.8423
$ws.Range("K74").Value = 4017.6365
$ws.Range("L74").Value = 4326.8423
$ws.Range("M74").Value = -3143.6365
$ws.Range("N74").Value = -6074.8423

# Sheet ARM, row 77 (G77=44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4160.927
$ws.Range("I77").Value = 4017.6365
$ws.Range("J77").Value = 4326.8423
$ws.Range("K77").Value = 20088.1825
$ws.Range("L77").Value = 21634.2115
$ws.Range("M77").Value = -15720.1825
$ws.Range("N77").Value = -30370.2115

# Sheet ARM, row 116 (G116=27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 6864.9
$ws.Range("I116").Value = 5223.5884
$ws.Range("K116").Value = 5223.5884
$ws.Range("M116").Value = -2929.5884

# Sheet BSM, row 3 (G3=27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6864.9
$ws.Range("I3").Value = 5223.5884
$ws.Range("K3").Value = 5223.5884
$ws.Range("M3").Value = -5109.5884

# Sheet BSM, row 23 (G23=1686)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 5950
$ws.Range("I23").Value = 5900
$ws.Range("K23").Value = 5900
$ws.Range("M23").Value = -5617

# Sheet BSM, row 86 (G86=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 325599
$ws.Range("I86").Value = 558563
$ws.Range("J86").Value = 3033.4614
$ws.Range("K86").Value = 558563
$ws.Range("L86").Value = 3033.4614
$ws.Range("M86").Value = -557440
$ws.Range("N86").Value = -5279.4614

# Sheet BSM, row 89 (G89=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 325599
$ws.Range("I89").Value = 558563
$ws.Range("J89").Value = 3033.4614
$ws.Range("K89").Value = 2792815
$ws.Range("L89").Value = 15167.307
$ws.Range("M89").Value = -2787199
$ws.Range("N89").Value = -26399.307

# Sheet BSM, row 94 (G94=19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1660.1666
$ws.Range("I94").Value = 1679.909
$ws.Range("J94").Value = 1443
$ws.Range("K94").Value = 1679.909
$ws.Range("L94").Value = 1443
$ws.Range("M94").Value = -1228.909
$ws.Range("N94").Value = -2345

# Sheet BSM, row 134 (G134=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7602.359
$ws.Range("I134").Value = 4932.2354
$ws.Range("J134").Value = 25759.2
$ws.Range("K134").Value = 14796.7062
$ws.Range("L134").Value = 77277.60000000001
$ws.Range("M134").Value = -12261.7062
$ws.Range("N134").Value = -82347.60000000001

# Sheet CRP, row 38 (G38=1637)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Sheet CRP, row 39 (G39=1915)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 4597
$ws.Range("I39").Value = 4597
$ws.Range("K39").Value = 4597
$ws.Range("M39").Value = -4206

# Sheet CRP, row 46 (G46=1637)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Sheet CRP, row 49 (G49=1915)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 4597
$ws.Range("I49").Value = 4597
$ws.Range("K49").Value = 4597
$ws.Range("M49").Value = -4415

# Sheet CRP, row 86 (G86=12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2825.842
$ws.Range("I86").Value = 2750.5
$ws.Range("J86").Value = 2880.6365
$ws.Range("K86").Value = 2750.5
$ws.Range("L86").Value = 2880.6365
$ws.Range("M86").Value = -1627.5
$ws.Range("N86").Value = -5126.636500000001

# Sheet CRP, row 89 (G89=12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2825.842
$ws.Range("I89").Value = 2750.5
$ws.Range("J89").Value = 2880.6365
$ws.Range("K89").Value = 13752.5
$ws.Range("L89").Value = 14403.1825
$ws.Range("M89").Value = -8136.5
$ws.Range("N89").Value = -25635.1825

# Sheet CRP, row 122 (G122=36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1420.7778
$ws.Range("I122").Value = 1231.1666
$ws.Range("K122").Value = 3693.4998
$ws.Range("M122").Value = -1243.4998

# Sheet CUL, row 22 (G22=4697)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 125
$ws.Range("J22").Value = 150
$ws.Range("L22").Value = 450
$ws.Range("N22").Value = -788

# Sheet CUL, row 27 (G27=4697)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 125
$ws.Range("J27").Value = 150
$ws.Range("L27").Value = 450
$ws.Range("N27").Value = -654

# Sheet CUL, row 127 (G127=38263)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3598.5
$ws.Range("J127").Value = 3598.5
$ws.Range("L127").Value = 10795.5
$ws.Range("N127").Value = -20715.5

# Sheet CUL, row 132 (G132=43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 11114855
$ws.Range("I132").Value = 1938.4
$ws.Range("K132").Value = 17445.6
$ws.Range("M132").Value = -14915.6

# Sheet GSM, row 12 (G12=4093)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 8589
$ws.Range("I12").Value = 7649.3335
$ws.Range("J12").Value = 9998.5
$ws.Range("K12").Value = 7649.3335
$ws.Range("L12").Value = 9998.5
$ws.Range("M12").Value = -7509.3335
$ws.Range("N12").Value = -10278.5

# Sheet GSM, row 33 (G33=4450)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 29446.25
$ws.Range("J33").Value = 29446.25
$ws.Range("L33").Value = 29446.25
$ws.Range("N33").Value = -29950.25

# Sheet GSM, row 35 (G35=4317)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 27814.857
$ws.Range("I35").Value = 18724.75
$ws.Range("J35").Value = 39935
$ws.Range("K35").Value = 18724.75
$ws.Range("L35").Value = 39935
$ws.Range("M35").Value = -18426.75
$ws.Range("N35").Value = -40531

# Sheet GSM, row 36 (G36=4222)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 11649.333
$ws.Range("J36").Value = 11649.333
$ws.Range("L36").Value = 11649.333
$ws.Range("N36").Value = -12619.333

# Sheet GSM, row 80 (G80=12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10087.25
$ws.Range("I80").Value = 5749.5
$ws.Range("J80").Value = 14425
$ws.Range("K80").Value = 5749.5
$ws.Range("L80").Value = 14425
$ws.Range("M80").Value = -4751.5
$ws.Range("N80").Value = -16421

# Sheet GSM, row 83 (G83=12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10087.25
$ws.Range("I83").Value = 5749.5
$ws.Range("J83").Value = 14425
$ws.Range("K83").Value = 28747.5
$ws.Range("L83").Value = 72125
$ws.Range("M83").Value = -23755.5
$ws.Range("N83").Value = -82109

# Sheet GSM, row 132 (G132=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 22599.092
$ws.Range("I132").Value = 41524.6
$ws.Range("J132").Value = 6827.8335
$ws.Range("K132").Value = 124573.8
$ws.Range("L132").Value = 20483.5005
$ws.Range("M132").Value = -122043.8
$ws.Range("N132").Value = -25543.5005

# Sheet LTW, row 7 (G7=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12197.733
$ws.Range("I7").Value = 14998.454
$ws.Range("K7").Value = 14998.454
$ws.Range("M7").Value = -14886.454

# Sheet LTW, row 41 (G41=3611)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 8499.5
$ws.Range("J41").Value = 8499.5
$ws.Range("L41").Value = 8499.5
$ws.Range("N41").Value = -9375.5

# Sheet LTW, row 68 (G68=12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4699.6665
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 6059.4
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 6059.4
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -7557.4

# Sheet LTW, row 71 (G71=12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4699.6665
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 6059.4
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 30297
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -37785

# Sheet LTW, row 100 (G100=19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5700.375
$ws.Range("I100").Value = 6081.231
$ws.Range("J100").Value = 4050
$ws.Range("K100").Value = 6081.231
$ws.Range("L100").Value = 4050
$ws.Range("M100").Value = -5540.231
$ws.Range("N100").Value = -5132

# Sheet LTW, row 126 (G126=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 12197.733
$ws.Range("I126").Value = 14998.454
$ws.Range("K126").Value = 44995.362
$ws.Range("M126").Value = -42525.362

# Sheet LTW, row 132 (G132=44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4845.5483
$ws.Range("I132").Value = 4154.75
$ws.Range("K132").Value = 12464.25
$ws.Range("M132").Value = -9934.25

# Sheet LTW, row 136 (G136=44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4009.532
$ws.Range("I136").Value = 3539.8572
$ws.Range("K136").Value = 10619.5716
$ws.Range("M136").Value = -8069.571599999999

# Sheet WVR, row 17 (G17=3539)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4056.4443
$ws.Range("I17").Value = 3784.6667
$ws.Range("J17").Value = 4600
$ws.Range("K17").Value = 3784.6667
$ws.Range("L17").Value = 4600
$ws.Range("M17").Value = -3612.6667
$ws.Range("N17").Value = -4944

# Sheet WVR, row 20 (G20=3023)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

# Sheet WVR, row 21 (G21=3341)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 22866.428
$ws.Range("J21").Value = 22866.428
$ws.Range("L21").Value = 22866.428
$ws.Range("N21").Value = -23336.428

# Sheet WVR, row 23 (G23=3325)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2216.8333
$ws.Range("I23").Value = 50.25
$ws.Range("J23").Value = 6550
$ws.Range("K23").Value = 50.25
$ws.Range("L23").Value = 6550
$ws.Range("M23").Value = 178.75
$ws.Range("N23").Value = -7008

# Sheet WVR, row 35 (G35=3341)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 22866.428
$ws.Range("J35").Value = 22866.428
$ws.Range("L35").Value = 22866.428
$ws.Range("N35").Value = -23446.428

# Sheet WVR, row 42 (G42=3372)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 66260.25
$ws.Range("J42").Value = 85021.5
$ws.Range("L42").Value = 85021.5
$ws.Range("N42").Value = -85777.5

# Sheet WVR, row 132 (G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1304759.6
$ws.Range("I132").Value = 1742535.9
$ws.Range("J132").Value = 31228.773
$ws.Range("K132").Value = 5227607.699999999
$ws.Range("L132").Value = 93686.319
$ws.Range("M132").Value = -5225077.699999999
$ws.Range("N132").Value = -98746.319
